# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column D) for the last file row (row 5,
# d33ba4ab-...) on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-08 07:46:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-08 07:46:28"
